$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.665.57'
$ws.Range('E2').Value = '  -4.13%  '

$ws.Range('D3').Value = '2.568.98'
$ws.Range('E3').Value = '  -3.54%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '511.82'
$ws.Range('E5').Value = '  -4.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.30'
$ws.Range('E6').Value = '  -6.16%  '

$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.575'
$ws.Range('E8').Value = '  -3.03%  '

$ws.Range('D9').Value = '2.584.84'
$ws.Range('E9').Value = '  -3.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.37'
$ws.Range('E10').Value = '  -3.77%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  -5.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.337'
$ws.Range('E12').Value = '  -4.63%  '

$ws.Range('E13').Value = '  -0.83%  '

$ws.Range('D14').Value = '3.020.63'
$ws.Range('E14').Value = '  -3.42%  '

$ws.Range('D15').Value = '58.597.21'
$ws.Range('E15').Value = '  -4.17%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.30'
$ws.Range('E16').Value = '  -3.96%  '

$ws.Range('E17').Value = '  -4.61%  '

$ws.Range('D18').Value = '2.576.57'
$ws.Range('E18').Value = '  -3.50%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '348.88'
$ws.Range('E19').Value = '  -1.99%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.57'
$ws.Range('E20').Value = '  -4.50%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.35'
$ws.Range('E21').Value = '  -3.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.06'
$ws.Range('E22').Value = '  -3.73%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.72'
$ws.Range('E24').Value = '  -1.52%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.419'
$ws.Range('E25').Value = '  -3.38%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.36%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -4.43%  '

$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.680.31'
$ws.Range('E28').Value = '  -3.52%  '

$ws.Range('D29').Value = '0.0₃0817'
$ws.Range('E29').Value = '  -5.57%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.05'
$ws.Range('E30').Value = '  -5.16%  '

$ws.Range('E31').Value = '  -0.02%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.00'
$ws.Range('E32').Value = '  -3.44%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.77'
$ws.Range('E33').Value = '  -4.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.51'
$ws.Range('E34').Value = '  -0.75%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').Value = '  -5.28%  '

$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.933'
$ws.Range('E36').Value = '  +5.42%  '

$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.98'
$ws.Range('E37').Value = '  -4.40%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.14'
$ws.Range('E38').Value = '  -5.43%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.857'
$ws.Range('E39').Value = '  -7.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.11'
$ws.Range('E40').Value = '  -1.90%  '

$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '292.03'
$ws.Range('E41').Value = '  -5.32%  '

$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.42'
$ws.Range('E42').Value = '  -5.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.59'
$ws.Range('E43').Value = '  -6.33%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0999'
$ws.Range('E44').Value = '  -2.70%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.996'
$ws.Range('E45').Value = '  -0.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.610'
$ws.Range('E46').Value = '  -6.44%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0540'
$ws.Range('E47').Value = '  -4.94%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.30'
$ws.Range('E48').Value = '  -5.05%  '

$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0231'
$ws.Range('E49').Value = '  -3.94%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.71'
$ws.Range('E50').Value = '  -6.22%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.25'
$ws.Range('E51').Value = '  -1.03%  '
